$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Transmitance values (column B, rows 3-18) from 1 to 100
$ws.Range("B3:B18").Value = 100

# Update the selected/active cell to B18 (matches the saved view state in the diff)
$ws.Range("B18").Select()
